$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date updated
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value now populated
$meta.Range("B9").Value = "Alvearie Team"

# The sheet had a duplicated "Contact / No display for ContactDetail" row
# (rows 10 and 11). Remove the first duplicate row entirely, which shifts
# everything below it up by one.
$meta.Rows.Item(10).Delete()

# The remaining former-duplicate row (now row 10) becomes the new
# "Jurisdiction / United States of America" row.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short & Definition now reflect the profile itself
# instead of the generic placeholder text.
$elements.Range("K2").Value = "Sex Assigned At Birth"
$elements.Range("L2").Value = "Sex assigned at birth based on observation by a physician.  Also called phenotype, or how the genes were expressed, based on the observation of the doctor. This assignment may not match the gender identity of the person, which will not be known until the newborn is older."
